# Adicionando check botao concluir
# Append a new "objeto" row (ar sala / A/C, 23, concluido=FALSE) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ar sala"
$ws.Range("B3").Value = "A/C"
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = $false
